$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data row (2-51) had its 'Hora' (G) value bumped from 20 to 21 on refresh.
$ws.Range("G2:G51").Value = '''21'

# Row 2
$ws.Range("D2").Value = '''246.55'

# Row 3
$ws.Range("D3").Value = '''22.37'

# Row 6
$ws.Range("D6").Value = '''3.413'

# Row 7
$ws.Range("D7").Value = '''6.308'

# Row 8
$ws.Range("D8").Value = '''0.8061'

# Row 9
$ws.Range("D9").Value = '''0.8575'

# Row 10
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '''0.01094'
$ws.Range("E10").Value = '9OneONE'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1412'
$ws.Range("E11").Value = '10WazirXWRX'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07349'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '''0.03026'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '''0.03076'
$ws.Range("E14").Value = '13BitrueCoinBTR'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '''0.09387'
$ws.Range("E15").Value = '14BitMartTokenBMX'

# Row 16
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '''3.880'
$ws.Range("E16").Value = '15MCDexMCB'

# Row 17
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '''0.001580'
$ws.Range("E17").Value = '16BitForexTokenBF'

# Row 18
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '''0.04777'
$ws.Range("E18").Value = '17CoinExTokenCET'

# Row 19
$ws.Range("D19").Value = '''0.006408'

# Row 20
$ws.Range("D20").Value = '''0.005025'

# Row 22
$ws.Range("D22").Value = '''0.0001500'

# Row 23
$ws.Range("D23").Value = '''3.694'

# Row 24
$ws.Range("D24").Value = '''2.192'

# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006820'
$ws.Range("E41").Value = '40KickTokenKICK'

# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1065'
$ws.Range("E42").Value = '41BKEXTokenBKK'

# Row 43
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.003200'
$ws.Range("E43").Value = '42CEJICEJI'

# Row 44
$ws.Range("D44").Value = '''0.008468'

# Row 45
$ws.Range("D45").Value = '''0.00005594'

# Row 48
$ws.Range("D48").Value = '''0.1964'
$ws.Range("E48").Value = '47BOLOBOLOWorstin24h'
